# Add a new "ScheduleName" / "AstroSchedule" row to the Configuration sheet
# and make that sheet the active tab (matching the author's edit).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Configuration")

# New row of configuration data, directly below the existing EBOVersion row.
$ws.Range("A7").Value = "ScheduleName"
$ws.Range("B7").Value = "AstroSchedule"

# Column B needs to widen to fit the new, longer value.
$ws.Columns.Item(2).AutoFit()

# The Configuration sheet becomes the active/visible tab ...
$ws.Activate()

# ... with the selection left just below the newly entered data.
$ws.Range("B8").Select() | Out-Null
